$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price rows for "Comercializadora del Agro de Limarí - Ciruela"
# appended as rows 12-14, mirroring the existing layout (Region de O'Higgins,
# week of 2022-02-03, $/caja 15 kilos granel).

$rows = @(
    @{ Row = 12; Calidad = "Especial"; Volumen = 160; PrecioMin = 15500; PrecioMax = 16000; PrecioProm = 15750; PrecioKg = 1050 },
    @{ Row = 13; Calidad = "Primera";  Volumen = 200; PrecioMin = 13500; PrecioMax = 14000; PrecioProm = 13750; PrecioKg = 917 },
    @{ Row = 14; Calidad = "Segunda";  Volumen = 300; PrecioMin = 10500; PrecioMax = 11000; PrecioProm = 10750; PrecioKg = 717 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Cells.Item($i, 1).Value = 2
    $ws.Cells.Item($i, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($i, 3).Value = "Coquimbo"

    $ws.Cells.Item($i, 4).Value = 44595
    $ws.Cells.Item($i, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($i, 5).Value = 4
    $ws.Cells.Item($i, 6).Value = "Fruta"
    $ws.Cells.Item($i, 7).Value = 100103
    $ws.Cells.Item($i, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($i, 9).Value = 100103002
    $ws.Cells.Item($i, 10).Value = "Ciruela"
    $ws.Cells.Item($i, 11).Value = "Black Amber"
    $ws.Cells.Item($i, 12).Value = $r.Calidad
    $ws.Cells.Item($i, 13).Value = $r.Volumen
    $ws.Cells.Item($i, 14).Value = $r.PrecioMin
    $ws.Cells.Item($i, 15).Value = $r.PrecioMax
    $ws.Cells.Item($i, 16).Value = $r.PrecioProm
    $ws.Cells.Item($i, 17).Value = "`$/caja 15 kilos granel"
    $ws.Cells.Item($i, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($i, 19).Value = $r.PrecioKg
    $ws.Cells.Item($i, 20).Value = 15
}
